# Updated symbol list - refresh Price / Volume(1h) columns with latest quotes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for the Price (D) and Volume(1h) (E) columns, keyed by cell address.
# These are stored in the sheet as literal text (e.g. "273.69", "-1.85%"),
# so each write keeps the cell formatted as Text (NumberFormat "@") to stop Excel
# from auto-coercing the numeric-looking / percent-looking strings into numbers,
# then restores the "Normal" style so no stray formatting is left behind.
$updates = [ordered]@{
    "D2" = "273.69"
    "E2" = "-1.85%"
    "D3" = "26.65"
    "E3" = "-2.25%"
    "D4" = "4.887"
    "E4" = "1.84%"
    "D5" = "0.06322"
    "E5" = "1.12%"
    "D7" = "3.360"
    "E7" = "2.85%"
    "D8" = "1.252"
    "E8" = "32.63%"
    "D9" = "0.8766"
    "E9" = "-0.10%"
    "E10" = "0.72%"
    "D11" = "0.05131"
    "E11" = "1.01%"
    "D12" = "0.07309"
    "E12" = "0.33%"
    "D13" = "0.03123"
    "E13" = "-0.83%"
    "E14" = "0.07%"
    "D15" = "0.001575"
    "E15" = "1.64%"
    "D16" = "0.0006340"
    "E16" = "1.21%"
    "D17" = "0.006038"
    "E17" = "1.76%"
    "E18" = "-0.38%"
    "D19" = "2.284"
    "E19" = "0.39%"
    "E21" = "1.36%"
    "D22" = "3.908"
    "E22" = "1.58%"
    "D23" = "0.04413"
    "E23" = "2.02%"
    "D24" = "0.001177"
    "E24" = "0.19%"
    "D25" = "0.004404"
    "E25" = "3.06%"
    "D26" = "0.0001200"
    "E26" = "0.18%"
    "D27" = "0.0001704"
    "E27" = "5.65%"
    "D40" = "0.04021"
    "E40" = "-0.22%"
    "D41" = "0.006653"
    "E41" = "-0.80%"
    "E42" = "1.15%"
    "D43" = "0.002091"
    "E43" = "-5.26%"
    "D44" = "0.01225"
    "E44" = "-2.35%"
    "D45" = "0.00005312"
    "E45" = "3.70%"
    "D46" = "2.357"
    "E46" = "6.31%"
    "D47" = "0.02004"
    "E47" = "-32.91%"
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.Style = "Normal"
}
